# The sheet stores quarterly series as groups of 4 rows per year: A, B, C, D.
# In this edit the "B" and "C" rows of every year-group swap places (the
# "产销率"/"销售量" helper columns F and G - which were derived/redundant
# columns - are also dropped, shrinking the used range from A1:G81 to A1:E81).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($yearStart = 2; $yearStart -le 78; $yearStart += 4) {
    $rowB = $yearStart + 1
    $rowC = $yearStart + 2

    $rangeB = "A" + $rowB + ":E" + $rowB
    $rangeC = "A" + $rowC + ":E" + $rowC

    $tmp = $ws.Range($rangeB).Value2
    $ws.Range($rangeB).Value2 = $ws.Range($rangeC).Value2
    $ws.Range($rangeC).Value2 = $tmp
}

# Drop the trailing 产销率 / 销售量 columns (F:G) entirely.
$ws.Range("F1:G81").Delete()
